$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the "Vehicle Registration State" row (row 46) to
# hold the new "Vehicle Registration Non-Expiring Indicator" mapping entry.
$ws.Rows("46").Insert()

# Populate the new row's CLASS/ELEMENT NAME and NIEM mapping cells.
$ws.Range("C46").Value = "Vehicle Registration Non-Expiring Indicator"
$ws.Range("F46").Value = "/wm-req-doc:WarrantIssuedReport/j:ConveyanceRegistration[not(j:RegistrationExpirationDate)]/wm-req-ext:ConveyanceRegistrationNonExpiringIndicator"

# The inserted row inherited a blank, styled "Sample Data" (column E) cell
# from the row above it; this row has no sample data, so drop it entirely.
$ws.Range("E46").Clear()

# Match the row height used by the surrounding mapping rows.
$ws.Rows("46").RowHeight = 56

# Leave the selection where Excel would land after entering the new row.
$ws.Range("F47").Select()
